$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.498.49'
$ws.Range('E2').Value = '  -2.13%  '

$ws.Range('D3').Value = '1.747.56'
$ws.Range('E3').Value = '  -3.20%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '324.04'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.04%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.12%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4437'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +3.42%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3606'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -0.79%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '42.28'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.73%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07460'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.53%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.096'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.03%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.002'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.26%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.59'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -5.09%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.027'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.54%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.135'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.63%  '

$ws.Range('D16').Value = '1.748.39'
$ws.Range('E16').Value = '  -4.43%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '92.19'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.15%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001058'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.05%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06395'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.17%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.11%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '16.82'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.65%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.848'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -5.15%  '

$ws.Range('D23').Value = '27.566.42'
$ws.Range('E23').Value = '  -2.04%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.18'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -2.67%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.091'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.26%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '161.60'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.79%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.38'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.58%  '

$ws.Range('D28').Value = '1.958.47'
$ws.Range('E28').Value = '  -3.60%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.103'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -6.05%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.63'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.88%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.074'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -8.77%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.657'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.36%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.08985'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.06%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.510'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -6.99%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '11.95'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -6.81%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02317'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.86%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2086'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.29%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.6325'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.87%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05972'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.64%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '4.934'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.11%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.209'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.70%  '

$ws.Range('E42').Value = '  +0.19%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.390'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.47%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '7.753'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.84%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.22'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -4.01%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.712'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.27%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5864'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.77%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '121.08'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.76%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.943'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.60%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.152'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.83%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06863'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.82%  '
